# Adds a new "CC_exp" (capacity-credit) row to the assumptions table on every
# scenario sheet, inserted just above the existing "C_exp" row (i.e. as the
# new row 4), pushing the rows that used to be 4-9 down to 5-10.
#
# Column layout for the new row: A=CC_exp, B=[-], C=<scenario value>,
# D=constant, E=0, F=0, G=0 (H left blank, matching the "script overwrites"
# H note shifting down with whatever row it was attached to).

$wb = $excel.ActiveWorkbook

# C-column value for the new "CC_exp" row, per sheet (1-based sheet index).
$values = @{
    1 = 0     # wind_only
    2 = 0.4   # 4_hr_batt
    3 = 1     # 10_hr_batt
    4 = 1     # 10_hr_ocaes
    5 = 1     # 24_hr_ocaes
    6 = 1     # 48_hr_ocaes
    7 = 1     # 72_hr_ocaes
    8 = 1     # 168_hr_ocaes
}

# Whether the new row's C cell should carry the "highlighted" (user-entered)
# formatting that the other populated rows in that column use. Only
# wind_only keeps the plain/default formatting (its whole C column is 0 and
# unhighlighted in the source file).
$highlight = @{
    1 = $false
    2 = $true
    3 = $true
    4 = $true
    5 = $true
    6 = $true
    7 = $true
    8 = $true
}

for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Insert a new blank row above the current row 4 (shifts old rows 4-9 to 5-10).
    $ws.Rows.Item(4).Insert()

    if ($highlight[$i]) {
        # Pick up the "highlighted" number format/fill used by the other
        # populated cells in column C (e.g. the row that is now C8 -
        # originally row 7, "L_well" - which always carries that formatting
        # already on every non-wind_only sheet) without disturbing its value.
        $ws.Range("C8").Copy()
        $ws.Range("C4").PasteSpecial(-4122)
        $excel.CutCopyMode = $false
    }

    $ws.Range("A4").Value2 = "CC_exp"
    $ws.Range("B4").Value2 = "[-]"
    $ws.Range("C4").Value2 = $values[$i]
    $ws.Range("D4").Value2 = "constant"
    $ws.Range("E4").Value2 = 0
    $ws.Range("F4").Value2 = 0
    $ws.Range("G4").Value2 = 0
}
